# Fruta / hortaliza, semanal
# Weekly refresh of the "Hortaliza, Mapocho Venta Directa de Santiago -
# Zapallo italiano" price series: updates the Fecha, Volumen, Precio
# mínimo/máximo/promedio ponderado, Unidad de comercialización, Origen and
# Precio $/Kg (and Kg o Unidades) figures for the affected weekly rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44333
$ws.Range("J2").Value = 25
$ws.Range("M2").Value = 10400
$ws.Range("O2").Value = 'Provincia de Limarí'

# Row 3
$ws.Range("D3").Value = 44405
$ws.Range("J3").Value = 45
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 9000
$ws.Range("N3").Value = '$/caja 50 unidades'
$ws.Range("O3").Value = 'Provincia de Quillota'
$ws.Range("P3").Value = 180
$ws.Range("Q3").Value = 50

# Row 4
$ws.Range("D4").Value = 44291
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 9000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 9000
$ws.Range("P4").Value = 150

# Row 5
$ws.Range("D5").Value = 44312
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("N5").Value = '$/caja 60 unidades'
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 167
$ws.Range("Q5").Value = 60

# Row 6
$ws.Range("D6").Value = 44186
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("P6").Value = 117

# Row 7
$ws.Range("D7").Value = 44179
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("P7").Value = 117

# Row 8
$ws.Range("D8").Value = 44243
$ws.Range("J8").Value = 80
$ws.Range("M8").Value = 10375
$ws.Range("O8").Value = 'Provincia de Quillota'

# Row 9
$ws.Range("D9").Value = 44277
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("P9").Value = 167

# Row 11
$ws.Range("D11").Value = 44315
$ws.Range("J11").Value = 25

# Row 12
$ws.Range("D12").Value = 44200
$ws.Range("J12").Value = 10
$ws.Range("K12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 9000
$ws.Range("P12").Value = 150
